$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Income Statement, Balance Sheet section references, and Cash Flow Statement
# with the latest quarterly column inserted (shift right) and refreshed figures.

# Row 8
$ws.Range("D8").Value = 20700
$ws.Range("E8").Value = 25200
$ws.Range("F8").Value = 26400
$ws.Range("G8").Value = 48300
$ws.Range("H8").Value = 13600
$ws.Range("I8").Value = 17300
$ws.Range("J8").Value = 26300

# Row 9
$ws.Range("D9").Value = 500
$ws.Range("E9").Value = 1700
$ws.Range("F9").Value = 2700
$ws.Range("G9").Value = 1900
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 100
$ws.Range("J9").Value = 100

# Row 10
$ws.Range("D10").Value = 20200
$ws.Range("E10").Value = 23500
$ws.Range("F10").Value = 23700
$ws.Range("G10").Value = 46400
$ws.Range("H10").Value = 13500
$ws.Range("I10").Value = 17200
$ws.Range("J10").Value = 26200

# Row 12
$ws.Range("D12").Value = 84900
$ws.Range("E12").Value = 67000
$ws.Range("F12").Value = 61000
$ws.Range("G12").Value = 72800
$ws.Range("H12").Value = 51800
$ws.Range("I12").Value = 84700
$ws.Range("J12").Value = 100600

# Row 14
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1800

# Row 17
$ws.Range("D17").Value = 114900
$ws.Range("E17").Value = 82600
$ws.Range("F17").Value = 75000
$ws.Range("G17").Value = 87700
$ws.Range("H17").Value = 59400
$ws.Range("I17").Value = 48000
$ws.Range("J17").Value = 60200

# Row 18
$ws.Range("D18").Value = -94200
$ws.Range("E18").Value = -57400
$ws.Range("F18").Value = -48600
$ws.Range("G18").Value = -39400
$ws.Range("H18").Value = -45800
$ws.Range("I18").Value = -30700
$ws.Range("J18").Value = -33900

# Row 20
$ws.Range("D20").Value = 2100
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = -100
$ws.Range("G20").Value = 1600
$ws.Range("H20").Value = 5600
$ws.Range("I20").Value = 5400
$ws.Range("J20").Value = -500

# Row 21
$ws.Range("D21").Value = -91400
$ws.Range("E21").Value = -56700
$ws.Range("F21").Value = -48500
$ws.Range("G21").Value = -37600
$ws.Range("H21").Value = -39900
$ws.Range("I21").Value = -24700
$ws.Range("J21").Value = "NA"

# Row 22
$ws.Range("D22").Value = 9200
$ws.Range("E22").Value = 8600
$ws.Range("F22").Value = 6500
$ws.Range("G22").Value = 5200
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 4800
$ws.Range("J22").Value = 4700

# Row 23
$ws.Range("D23").Value = -101300
$ws.Range("E23").Value = -65800
$ws.Range("F23").Value = -55100
$ws.Range("G23").Value = -43000
$ws.Range("H23").Value = -45200
$ws.Range("I23").Value = -30100
$ws.Range("J23").Value = -39100

# Row 26
$ws.Range("D26").Value = -101300
$ws.Range("E26").Value = -65800
$ws.Range("F26").Value = -55100
$ws.Range("G26").Value = -43000
$ws.Range("H26").Value = -45200
$ws.Range("I26").Value = -30100
$ws.Range("J26").Value = -39100

# Row 27
$ws.Range("D27").Value = -101300
$ws.Range("E27").Value = -65800
$ws.Range("F27").Value = -55100
$ws.Range("G27").Value = -43000
$ws.Range("H27").Value = -45200
$ws.Range("I27").Value = -30100
$ws.Range("J27").Value = -39100

# Row 32
$ws.Range("D32").Value = -2100
$ws.Range("E32").Value = -200
$ws.Range("F32").Value = 100
$ws.Range("G32").Value = -1600
$ws.Range("H32").Value = -5600
$ws.Range("I32").Value = -5400
$ws.Range("J32").Value = 500

# Row 33
$ws.Range("D33").Value = -101300
$ws.Range("E33").Value = -65800
$ws.Range("F33").Value = -55100
$ws.Range("G33").Value = -43000
$ws.Range("H33").Value = -45200
$ws.Range("I33").Value = -30100
$ws.Range("J33").Value = -39100

# Row 35
$ws.Range("D35").Value = -101300
$ws.Range("E35").Value = -65800
$ws.Range("F35").Value = -55100
$ws.Range("G35").Value = -43000
$ws.Range("H35").Value = -45200
$ws.Range("I35").Value = -30100
$ws.Range("J35").Value = -39100

# Row 81
$ws.Range("D81").Value = -101300
$ws.Range("E81").Value = -65800
$ws.Range("F81").Value = -55100
$ws.Range("G81").Value = -43000
$ws.Range("H81").Value = -45200
$ws.Range("I81").Value = -30100
$ws.Range("J81").Value = -39100

# Row 83
$ws.Range("J83").Value = "NA"

# Row 94
$ws.Range("J94").Value = "NA"

# Row 100
$ws.Range("J100").Value = "NA"

# Row 101
$ws.Range("J101").Value = "NA"
